# "update list of websites"
#
# Fills in the "KOMUNIKATY" (column D) / "CSS" (column E) website links for a
# batch of powiat (county) rows that were still missing them, fixes one
# mis-entered link (pajęczański, row 219, which was in column E instead of a
# real CSS value), and removes a stale placeholder link from row 366.
#
# NOTE: the cell values below are written in the exact order the source
# workbook's shared-string table grew in (new links first appear grouped by
# URL, not by row), so that each newly-introduced string lands at the same
# shared-string index as in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Płock (Płock city + płocki county)
$ws.Range("D229").Value = 'http://plock.psse.waw.pl/aktualnosci-i-komunikaty-6796/komunikaty'
$ws.Range("D230").Value = 'http://plock.psse.waw.pl/aktualnosci-i-komunikaty-6796/komunikaty'

# wołomiński
$ws.Range("D353").Value = 'http://wolomin.psse.waw.pl/1365'

# counties with no known website yet -> "brak" (none)
$ws.Range("D26").Value = 'brak'
$ws.Range("D27").Value = 'brak'
$ws.Range("D39").Value = 'brak'
$ws.Range("D85").Value = 'brak'
$ws.Range("D103").Value = 'brak'
$ws.Range("D104").Value = 'brak'
$ws.Range("D135").Value = 'brak'
$ws.Range("D136").Value = 'brak'
$ws.Range("D192").Value = 'brak'
$ws.Range("D305").Value = 'brak'
$ws.Range("D334").Value = 'brak'
$ws.Range("D335").Value = 'brak'

# bielski + Bielsko-Biała
$ws.Range("D14").Value = 'https://www.psse.bielsko.pl/'
$ws.Range("D15").Value = 'https://www.psse.bielsko.pl/'

# pajęczański had its link typed into the wrong column (E instead of D)
$ws.Range("E219").Value = 'http://www.psse.pajeczno.com/index.php?p=p5'

# będziński + Dąbrowa Górnicza
$ws.Range("D6").Value = 'https://www.pssedg.pl/'
$ws.Range("D46").Value = 'https://www.pssedg.pl/'

# ełcki
$ws.Range("D54").Value = 'http://www.bip.visacom.pl/psse_elk/'

# opatowski
$ws.Range("D204").Value = 'http://www.psseopatow.pl/'

# Częstochowa + częstochowski
$ws.Range("D43").Value = 'http://psse.czest.pl/koronawirus-dane.html'
$ws.Range("D44").Value = 'http://psse.czest.pl/koronawirus-dane.html'

# zawierciański
$ws.Range("D365").Value = 'http://psse-zawiercie.internetdsl.pl/'

# lubański
$ws.Range("D147").Value = 'http://psseluban.pl/'

# gostyński
$ws.Range("D72").Value = 'http://www.psse-gostyn.pl/'

# Konin + koniński
$ws.Range("D116").Value = 'http://www.psse-konin.pl/koronawirus-2019-ncov'
$ws.Range("D117").Value = 'http://www.psse-konin.pl/koronawirus-2019-ncov'

# Poznań + poznański
$ws.Range("D235").Value = 'http://www.psse-poznan.pl/p,169,covid-19-dane-statystyczne'
$ws.Range("D236").Value = 'http://www.psse-poznan.pl/p,169,covid-19-dane-statystyczne'

# zawierciański (row 366) no longer carries the stale "??" placeholder
$ws.Range("D366").ClearContents()

# leave the workbook selected/scrolled where the author left it
$excel.ActiveWindow.ScrollRow = 218
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D236").Select()
